$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 67: "Convert Integers to the Sum of Two No-Zero Integers"
# (set first so new shared strings are appended in the same order Excel wrote them)
$ws.Range("A67").Value = 1317
$ws.Range("B67").Value = "Convert Integers to the Sum of Two No-Zero Integers"
$ws.Range("C67").Value = "Math"
$ws.Range("D67").Value = "Use a util function to check a zero in the digits/Loop until both number have no zero"

# Extend existing row 65 ("Fruits Into Baskets 2") with Type/Steps/extra steps
$ws.Range("C65").Value = "Iteration/Loop"
$ws.Range("D65").Value = "Loop in Loop"
$ws.Range("E65").Value = "Can use a track for usedBasket"

# Match the author's final selection/viewport state
$ws.Range("E65").Select()
